$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.866.04'
$ws.Range("E2").Value = '  +4.83%  '

$ws.Range("D3").Value = '''2.278.41'
$ws.Range("E3").Value = '  +2.22%  '

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").Value = '''231.35'
$ws.Range("E5").Value = '  -0.60%  '

$ws.Range("E6").Value = '  +1.52%  '

$ws.Range("D7").Value = '''64.13'
$ws.Range("E7").Value = '  +5.78%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").Value = '''0.426'
$ws.Range("E9").Value = '  +5.20%  '

$ws.Range("D10").Value = '''0.0978'
$ws.Range("E10").Value = '  +8.48%  '

$ws.Range("D11").Value = '''57.59'
$ws.Range("E11").Value = '  -1.07%  '

$ws.Range("D12").Value = '''26.22'
$ws.Range("E12").Value = '  +15.05%  '

$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("D14").Value = '''2.622.52'
$ws.Range("E14").Value = '  +2.48%  '

$ws.Range("D15").Value = '''15.80'
$ws.Range("E15").Value = '  +1.63%  '

$ws.Range("E16").Value = '  +4.96%  '

$ws.Range("D17").Value = '''0.820'
$ws.Range("E17").Value = '  +2.15%  '

$ws.Range("D18").Value = '''2.314.28'
$ws.Range("E18").Value = '  +3.04%  '

$ws.Range("D19").Value = '''43.704.28'
$ws.Range("E19").Value = '  +4.65%  '

$ws.Range("D20").Value = '''0.0₃0957'
$ws.Range("E20").Value = '  +5.19%  '

$ws.Range("D21").Value = '''73.44'
$ws.Range("E21").Value = '  +1.19%  '

$ws.Range("D22").Value = '''6.15'
$ws.Range("E22").Value = '  -0.44%  '

$ws.Range("D23").Value = '''250.62'
$ws.Range("E23").Value = '  +1.05%  '

$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").Value = '''2.53'
$ws.Range("E25").Value = '  +5.95%  '

$ws.Range("D26").Value = '''2.34'
$ws.Range("E26").Value = '  +1.05%  '

$ws.Range("D27").Value = '''9.96'
$ws.Range("E27").Value = '  +2.94%  '

$ws.Range("D28").Value = '''172.30'

$ws.Range("E29").Value = '  -1.67%  '

$ws.Range("E30").Value = '  +3.26%  '

$ws.Range("D31").Value = '''1.46'
$ws.Range("E31").Value = '  +3.87%  '

$ws.Range("D32").Value = '''2.79'
$ws.Range("E32").Value = '  +8.08%  '

$ws.Range("E33").Value = '  +0.42%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.0695'
$ws.Range("E34").Value = '  +5.57%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = '''5.16'
$ws.Range("E35").Value = '  +2.57%  '

$ws.Range("D36").Value = '''4.73'
$ws.Range("E36").Value = '  +0.44%  '

$ws.Range("D37").Value = '''6.85'
$ws.Range("E37").Value = '  +4.11%  '

$ws.Range("D38").Value = '''3.82'
$ws.Range("E38").Value = '  +5.55%  '

$ws.Range("D39").Value = '''2.35'
$ws.Range("E39").Value = '  -2.28%  '

$ws.Range("D40").Value = '''0.0249'
$ws.Range("E40").Value = '  +3.37%  '

$ws.Range("E41").Value = '  +0.31%  '

$ws.Range("D42").Value = '''10.94'
$ws.Range("E42").Value = '  +26.30%  '

$ws.Range("D43").Value = '''0.000225'
$ws.Range("E43").Value = '  -6.15%  '

$ws.Range("D44").Value = '''8.46'
$ws.Range("E44").Value = '  -2.05%  '

$ws.Range("D45").Value = '''4.58'
$ws.Range("E45").Value = '  +0.99%  '

$ws.Range("D46").Value = '''1.22'
$ws.Range("E46").Value = '  -0.35%  '

$ws.Range("D47").Value = '''0.0969'
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").Value = '''97.98'
$ws.Range("E48").Value = '  -1.08%  '

$ws.Range("D49").Value = '''1.487.15'
$ws.Range("E49").Value = '  +1.14%  '

$ws.Range("E50").Value = '  +1.69%  '

$ws.Range("D51").Value = '''2.32'
$ws.Range("E51").Value = '  +1.77%  '
